# Updates the cryptos list (prices + hourly volume %) as published by the
# "Updated cryptos list ... with GitHub Actions" workflow run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value to a cell while forcing it to stay text.
# Some "Price" values look numeric (e.g. "590.15", "0.440", "133.40") and
# Excel would otherwise silently convert them to a float (dropping
# trailing zeros / losing exact formatting) or reformat them. Temporarily
# switching the cell to the Text number format during assignment, then
# restoring its original style, keeps the text exact while leaving the
# cell's visual style untouched.
function Set-TextValue {
    param($row, $col, [string]$value)
    $cell = $ws.Cells.Item($row, $col)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $originalStyle
}

# Row 2 - Bitcoin
Set-TextValue 2 4 "60.925.23"
$ws.Cells.Item(2,5).Value = "  +0.11%  "

# Row 3 - Ethereum
Set-TextValue 3 4 "2.919.42"
$ws.Cells.Item(3,5).Value = "  +0.08%  "

# Row 4 - TetherUSD
$ws.Cells.Item(4,5).Value = "  -0.03%  "

# Row 5 - BNB
Set-TextValue 5 4 "590.15"
$ws.Cells.Item(5,5).Value = "  +1.03%  "

# Row 6 - Solana
Set-TextValue 6 4 "146.47"
$ws.Cells.Item(6,5).Value = "  +1.40%  "

# Row 7 - USDC
$ws.Cells.Item(7,5).Value = "  +0.01%  "

# Row 8 - XRP
$ws.Cells.Item(8,5).Value = "  +0.71%  "

# Row 9 - Toncoin
Set-TextValue 9 4 "6.88"
$ws.Cells.Item(9,5).Value = "  +0.87%  "

# Row 10 - Dogecoin
$ws.Cells.Item(10,5).Value = "  -0.69%  "

# Row 11 - Cardano
Set-TextValue 11 4 "0.440"
$ws.Cells.Item(11,5).Value = "  -1.46%  "

# Row 12 - ShibaInu
$ws.Cells.Item(12,5).Value = "  -0.28%  "

# Row 13 - Avalanche
$ws.Cells.Item(13,5).Value = "  -0.03%  "

# Row 14 - TRON
$ws.Cells.Item(14,5).Value = "  +0.10%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue 15 4 "3.402.51"
$ws.Cells.Item(15,5).Value = "  +0.08%  "

# Row 16 - WrappedBTC
Set-TextValue 16 4 "60.817.73"
$ws.Cells.Item(16,5).Value = "  +0.00%  "

# Row 17 - Polkadot
Set-TextValue 17 4 "6.70"
$ws.Cells.Item(17,5).Value = "  -0.70%  "

# Row 18 - WrappedEther
Set-TextValue 18 4 "2.918.52"
$ws.Cells.Item(18,5).Value = "  +0.01%  "

# Row 19 - BitcoinCash
Set-TextValue 19 4 "430.27"
$ws.Cells.Item(19,5).Value = "  -0.36%  "

# Row 20 - Chainlink
$ws.Cells.Item(20,5).Value = "  -1.96%  "

# Row 21 - Polygon
Set-TextValue 21 4 "0.679"
$ws.Cells.Item(21,5).Value = "  -0.68%  "

# Row 22 - Uniswap
Set-TextValue 22 4 "7.06"
$ws.Cells.Item(22,5).Value = "  -1.18%  "

# Row 23 - Litecoin
Set-TextValue 23 4 "81.41"
$ws.Cells.Item(23,5).Value = "  +1.23%  "

# Row 24 - RenderToken
$ws.Cells.Item(24,5).Value = "  +0.58%  "

# Row 25 - Fetch.AI
$ws.Cells.Item(25,5).Value = "  -0.90%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue 26 4 "11.86"
$ws.Cells.Item(26,5).Value = "  -0.49%  "

# Row 27 - Dai
$ws.Cells.Item(27,5).Value = "  +0.01%  "

# Row 28 - ImmutableX
Set-TextValue 28 4 "2.28"
$ws.Cells.Item(28,5).Value = "  +4.75%  "

# Row 29 - PancakeSwap
$ws.Cells.Item(29,5).Value = "  -0.07%  "

# Row 30 - NEARProtocol
$ws.Cells.Item(30,5).Value = "  -3.27%  "

# Row 31 - EthereumClassic
$ws.Cells.Item(31,5).Value = "  +0.29%  "

# Row 32 - Hedera
$ws.Cells.Item(32,5).Value = "  +2.25%  "

# Row 33 - FirstDigitalUSD
$ws.Cells.Item(33,5).Value = "  -0.04%  "

# Row 34 - PEPE
$ws.Cells.Item(34,5).Value = "  -1.01%  "

# Row 35 - Mantle
$ws.Cells.Item(35,5).Value = "  -0.04%  "

# Row 36 - Filecoin
$ws.Cells.Item(36,5).Value = "  -0.44%  "

# Row 37 - dogwifhat
$ws.Cells.Item(37,5).Value = "  +0.41%  "

# Row 38 - Stacks
$ws.Cells.Item(38,5).Value = "  -1.39%  "

# Row 39 - Kaspa
$ws.Cells.Item(39,5).Value = "  -3.74%  "

# Row 40 - Cosmos
$ws.Cells.Item(40,5).Value = "  -1.55%  "

# Row 41 - TheGraph
Set-TextValue 41 4 "0.282"
$ws.Cells.Item(41,5).Value = "  -4.64%  "

# Row 42 - Arweave
Set-TextValue 42 4 "40.05"
$ws.Cells.Item(42,5).Value = "  -2.89%  "

# Row 43 - Bittensor
Set-TextValue 43 4 "380.28"
$ws.Cells.Item(43,5).Value = "  +1.29%  "

# Row 44 - was Maker, now VeChain (rows 44/45 swapped content + new values)
$ws.Cells.Item(44,2).Value = "VeChain"
$ws.Cells.Item(44,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 44 4 "0.0343"
$ws.Cells.Item(44,5).Value = "  -1.59%  "

# Row 45 - was VeChain, now Maker
$ws.Cells.Item(45,2).Value = "Maker"
$ws.Cells.Item(45,3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 45 4 "2.692.84"
$ws.Cells.Item(45,5).Value = "  +0.72%  "

# Row 46 - Monero
Set-TextValue 46 4 "133.40"
$ws.Cells.Item(46,5).Value = "  +1.05%  "

# Row 48 - InjectiveProtocol
Set-TextValue 48 4 "23.80"
$ws.Cells.Item(48,5).Value = "  -2.97%  "

# Row 49 - Stellar
$ws.Cells.Item(49,5).Value = "  -0.62%  "

# Row 50 - ThetaToken
$ws.Cells.Item(50,5).Value = "  -3.59%  "

# Row 51 - Cronos
$ws.Cells.Item(51,5).Value = "  -0.18%  "
